$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02923533333333334
$ws.Range("H2").Value = 0.08770600000000001
$ws.Range("I2").Value = 0.01201286672177323
$ws.Range("J2").Value = 0.01201286672177324
$ws.Range("M2").Value = 92.253011
$ws.Range("N2").Value = 276.759033
$ws.Range("O2").Value = 0.2854710184133813
$ws.Range("P2").Value = 0.2854710184133813
$ws.Range("Q2").Value = 2.697047527588667
$ws.Range("R2").Value = 24.273427748298
$ws.Range("S2").Value = 0.003429325297128823
$ws.Range("T2").Value = 0.003429325297128823

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02923533333333334
$ws.Range("H3").Value = 0.08770600000000001
$ws.Range("I3").Value = 0.01201286672177323
$ws.Range("J3").Value = 0.01201286672177324
$ws.Range("O3").Value = 0.1080355352256351
$ws.Range("P3").Value = 0.1080355352256351
$ws.Range("Q3").Value = 1.020688456542667
$ws.Range("R3").Value = 9.186196108884001
$ws.Range("S3").Value = 0.001297816485880992
$ws.Range("T3").Value = 0.001297816485880992

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02923533333333334
$ws.Range("H4").Value = 0.08770600000000001
$ws.Range("I4").Value = 0.01201286672177323
$ws.Range("J4").Value = 0.01201286672177324
$ws.Range("M4").Value = 42.21774566666667
$ws.Range("N4").Value = 126.653237
$ws.Range("O4").Value = 0.1306401029076487
$ws.Range("P4").Value = 0.1306401029076487
$ws.Range("Q4").Value = 1.234249867146889
$ws.Range("R4").Value = 11.108248804322
$ws.Range("S4").Value = 0.001569362144748324
$ws.Range("T4").Value = 0.001569362144748324

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02923533333333334
$ws.Range("H5").Value = 0.08770600000000001
$ws.Range("I5").Value = 0.01201286672177323
$ws.Range("J5").Value = 0.01201286672177324
$ws.Range("M5").Value = 13.65158233333333
$ws.Range("N5").Value = 40.954747
$ws.Range("O5").Value = 0.0422439448794879
$ws.Range("P5").Value = 0.0422439448794879
$ws.Range("Q5").Value = 0.3991085600424444
$ws.Range("R5").Value = 3.591977040382
$ws.Range("S5").Value = 0.000507470879639223
$ws.Range("T5").Value = 0.000507470879639223

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02923533333333334
$ws.Range("H6").Value = 0.08770600000000001
$ws.Range("I6").Value = 0.01201286672177323
$ws.Range("J6").Value = 0.01201286672177324
$ws.Range("M6").Value = 18.17840666666667
$ws.Range("N6").Value = 54.53522
$ws.Range("O6").Value = 0.0562519120841046
$ws.Range("P6").Value = 0.05625191208410459
$ws.Range("Q6").Value = 0.5314517783688889
$ws.Range("R6").Value = 4.78306600532
$ws.Range("S6").Value = 0.0006757467227112538
$ws.Range("T6").Value = 0.0006757467227112538

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.02923533333333334
$ws.Range("H7").Value = 0.08770600000000001
$ws.Range("I7").Value = 0.01201286672177323
$ws.Range("J7").Value = 0.01201286672177324
$ws.Range("M7").Value = 121.9471053333333
$ws.Range("N7").Value = 365.841316
$ws.Range("O7").Value = 0.3773574864897424
$ws.Range("P7").Value = 0.3773574864897424
$ws.Range("Q7").Value = 3.565164273455112
$ws.Range("R7").Value = 32.086478461096
$ws.Range("S7").Value = 0.004533145191664619
$ws.Range("T7").Value = 0.00453314519166462

# Row 8
$ws.Range("H8").Value = 6.303437000000001
$ws.Range("I8").Value = 0.8633656599331188
$ws.Range("J8").Value = 0.8633656599331188
$ws.Range("M8").Value = 92.253011
$ws.Range("N8").Value = 276.759033
$ws.Range("O8").Value = 0.2854710184133813
$ws.Range("P8").Value = 0.2854710184133813
$ws.Range("Q8").Value = 193.8370142996024
$ws.Range("R8").Value = 1744.533128696421
$ws.Range("S8").Value = 0.2464658742042485
$ws.Range("T8").Value = 0.2464658742042484

# Row 9
$ws.Range("H9").Value = 6.303437000000001
$ws.Range("I9").Value = 0.8633656599331188
$ws.Range("J9").Value = 0.8633656599331188
$ws.Range("O9").Value = 0.1080355352256351
$ws.Range("P9").Value = 0.1080355352256351
$ws.Range("Q9").Value = 73.35695827473535
$ws.Range("R9").Value = 660.2126244726181
$ws.Range("S9").Value = 0.09327417116630814
$ws.Range("T9").Value = 0.09327417116630814

# Row 10
$ws.Range("H10").Value = 6.303437000000001
$ws.Range("I10").Value = 0.8633656599331188
$ws.Range("J10").Value = 0.8633656599331188
$ws.Range("M10").Value = 42.21774566666667
$ws.Range("N10").Value = 126.653237
$ws.Range("O10").Value = 0.1306401029076487
$ws.Range("P10").Value = 0.1306401029076487
$ws.Range("Q10").Value = 88.70563336395215
$ws.Range("R10").Value = 798.3507002755692
$ws.Range("S10").Value = 0.1127901786605927
$ws.Range("T10").Value = 0.1127901786605927

# Row 11
$ws.Range("H11").Value = 6.303437000000001
$ws.Range("I11").Value = 0.8633656599331188
$ws.Range("J11").Value = 0.8633656599331188
$ws.Range("M11").Value = 13.65158233333333
$ws.Range("N11").Value = 40.954747
$ws.Range("O11").Value = 0.0422439448794879
$ws.Range("P11").Value = 0.0422439448794879
$ws.Range("Q11").Value = 28.68396306282656
$ws.Range("R11").Value = 258.155667565439
$ws.Range("S11").Value = 0.03647197134905737
$ws.Range("T11").Value = 0.03647197134905736

# Row 12
$ws.Range("H12").Value = 6.303437000000001
$ws.Range("I12").Value = 0.8633656599331188
$ws.Range("J12").Value = 0.8633656599331188
$ws.Range("M12").Value = 18.17840666666667
$ws.Range("N12").Value = 54.53522
$ws.Range("O12").Value = 0.0562519120841046
$ws.Range("P12").Value = 0.05625191208410459
$ws.Range("Q12").Value = 38.19548039457112
$ws.Range("R12").Value = 343.75932355114
$ws.Range("S12").Value = 0.04856596919899275
$ws.Range("T12").Value = 0.04856596919899274

# Row 13
$ws.Range("H13").Value = 6.303437000000001
$ws.Range("I13").Value = 0.8633656599331188
$ws.Range("J13").Value = 0.8633656599331188
$ws.Range("M13").Value = 121.9471053333333
$ws.Range("N13").Value = 365.841316
$ws.Range("O13").Value = 0.3773574864897424
$ws.Range("P13").Value = 0.3773574864897424
$ws.Range("Q13").Value = 256.2286319336769
$ws.Range("R13").Value = 2306.057687403092
$ws.Range("S13").Value = 0.3257974953539194
$ws.Range("T13").Value = 0.3257974953539194

# Row 14
$ws.Range("G14").Value = 0.3032873333333334
$ws.Range("H14").Value = 0.9098620000000001
$ws.Range("I14").Value = 0.124621473345108
$ws.Range("J14").Value = 0.124621473345108
$ws.Range("M14").Value = 92.253011
$ws.Range("N14").Value = 276.759033
$ws.Range("O14").Value = 0.2854710184133813
$ws.Range("P14").Value = 0.2854710184133813
$ws.Range("Q14").Value = 27.97916969816067
$ws.Range("R14").Value = 251.812527283446
$ws.Range("S14").Value = 0.03557581891200402
$ws.Range("T14").Value = 0.03557581891200402

# Row 15
$ws.Range("G15").Value = 0.3032873333333334
$ws.Range("H15").Value = 0.9098620000000001
$ws.Range("I15").Value = 0.124621473345108
$ws.Range("J15").Value = 0.124621473345108
$ws.Range("O15").Value = 0.1080355352256351
$ws.Range("P15").Value = 0.1080355352256351
$ws.Range("Q15").Value = 10.58862153611867
$ws.Range("R15").Value = 95.29759382506802
$ws.Range("S15").Value = 0.01346354757344595
$ws.Range("T15").Value = 0.01346354757344595

# Row 16
$ws.Range("G16").Value = 0.3032873333333334
$ws.Range("H16").Value = 0.9098620000000001
$ws.Range("I16").Value = 0.124621473345108
$ws.Range("J16").Value = 0.124621473345108
$ws.Range("M16").Value = 42.21774566666667
$ws.Range("N16").Value = 126.653237
$ws.Range("O16").Value = 0.1306401029076487
$ws.Range("P16").Value = 0.1306401029076487
$ws.Range("Q16").Value = 12.80410750258823
$ws.Range("R16").Value = 115.236967523294
$ws.Range("S16").Value = 0.01628056210230771
$ws.Range("T16").Value = 0.01628056210230771

# Row 17
$ws.Range("G17").Value = 0.3032873333333334
$ws.Range("H17").Value = 0.9098620000000001
$ws.Range("I17").Value = 0.124621473345108
$ws.Range("J17").Value = 0.124621473345108
$ws.Range("M17").Value = 13.65158233333333
$ws.Range("N17").Value = 40.954747
$ws.Range("O17").Value = 0.0422439448794879
$ws.Range("P17").Value = 0.0422439448794879
$ws.Range("Q17").Value = 4.140352001657111
$ws.Range("R17").Value = 37.263168014914
$ws.Range("S17").Value = 0.005264502650791311
$ws.Range("T17").Value = 0.005264502650791311

# Row 18
$ws.Range("G18").Value = 0.3032873333333334
$ws.Range("H18").Value = 0.9098620000000001
$ws.Range("I18").Value = 0.124621473345108
$ws.Range("J18").Value = 0.124621473345108
$ws.Range("M18").Value = 18.17840666666667
$ws.Range("N18").Value = 54.53522
$ws.Range("O18").Value = 0.0562519120841046
$ws.Range("P18").Value = 0.05625191208410459
$ws.Range("Q18").Value = 5.513280482182223
$ws.Range("R18").Value = 49.61952433964001
$ws.Range("S18").Value = 0.007010196162400597
$ws.Range("T18").Value = 0.007010196162400597

# Row 19
$ws.Range("G19").Value = 0.3032873333333334
$ws.Range("H19").Value = 0.9098620000000001
$ws.Range("I19").Value = 0.124621473345108
$ws.Range("J19").Value = 0.124621473345108
$ws.Range("M19").Value = 121.9471053333333
$ws.Range("N19").Value = 365.841316
$ws.Range("O19").Value = 0.3773574864897424
$ws.Range("P19").Value = 0.3773574864897424
$ws.Range("Q19").Value = 36.98501238426578
$ws.Range("R19").Value = 332.8651114583921
$ws.Range("S19").Value = 0.04702684594415837
$ws.Range("T19").Value = 0.04702684594415837
